# Users.xlsx seed-data update ("Email Sending Service implemented.")
# Replaces the three test-user rows (email / password / role) with a new
# set of accounts and marks the touched cells with the Arial / theme-color-1
# font used elsewhere in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - admin account
$ws.Range("A1").Value = "pavlovt@gmail.com"
$ws.Range("B1").Value = "Test123!"
$ws.Range("C1").Value = "Admin"

# Row 2 - standard account
$ws.Range("A2").Value = "taskop@gmail.com"
$ws.Range("B2").Value = "Test123!"
$ws.Range("C2").Value = "Standard User"

# Row 3 - standard account
$ws.Range("A3").Value = "test@test.com"
$ws.Range("B3").Value = "Test123!"
$ws.Range("C3").Value = "Standard User"

# Re-apply the Arial / theme-color-1 font to the updated cells (matches the
# formatting already used in the workbook), row by row since this bridge
# does not support multi-area ("A1,A2:C2") range unions.
$ws.Range("B1").Font.Name = "Arial"
$ws.Range("B1").Font.ThemeColor = 1

$ws.Range("A2:C3").Font.Name = "Arial"
$ws.Range("A2:C3").Font.ThemeColor = 1
